# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.065.84"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Value = "'2.460.95"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'576.38"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "'146.48"
$ws.Range("E6").Value = "  +1.92%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.542"
$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("D9").Value = "'2.459.34"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "'5.29"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "'29.00"
$ws.Range("E14").Value = "  +9.25%  "

$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("D16").Value = "'2.906.36"
$ws.Range("E16").Value = "  +2.20%  "

$ws.Range("D17").Value = "'62.974.97"
$ws.Range("E17").Value = "  +2.76%  "

$ws.Range("D18").Value = "'2.457.37"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("E20").Value = "  +4.03%  "

$ws.Range("D21").Value = "'330.34"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("E22").Value = "  +12.45%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'66.53"
$ws.Range("E25").Value = "  +2.17%  "

$ws.Range("D26").Value = "'666.56"
$ws.Range("E26").Value = "  +7.95%  "

$ws.Range("D27").Value = "'1.16"
$ws.Range("E27").Value = "  +16.34%  "

$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  +8.01%  "

$ws.Range("D29").Value = "'0.0000100"
$ws.Range("E29").Value = "  +4.90%  "

$ws.Range("D31").Value = "'8.19"
$ws.Range("E31").Value = "  +1.47%  "

$ws.Range("D32").Value = "'1.45"
$ws.Range("E32").Value = "  +3.68%  "

$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +3.52%  "

$ws.Range("D34").Value = "'0.138"
$ws.Range("E34").Value = "  +3.77%  "

$ws.Range("E35").Value = "  +4.72%  "

$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "'4.80"
$ws.Range("E37").Value = "  +3.38%  "

$ws.Range("D38").Value = "'5.50"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").Value = "'153.15"
$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").Value = "'0.0₆0350"
$ws.Range("E42").Value = "  +23.97%  "

$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("D45").Value = "'42.31"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").Value = "'15.14"
$ws.Range("E47").Value = "  +27.55%  "

$ws.Range("D48").Value = "'146.72"
$ws.Range("E48").Value = "  +2.56%  "

$ws.Range("E49").Value = "  +2.18%  "

$ws.Range("D50").Value = "'20.78"
$ws.Range("E50").Value = "  +4.03%  "

$ws.Range("E51").Value = "  +1.84%  "
